# Insert a new data row before the existing row 208, pushing all the
# existing records (old rows 208-311) down by one row (new rows 209-312),
# and populate the new row 208 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(208).Insert()

$ws.Cells.Item(208, 1).Value2  = 5
$ws.Cells.Item(208, 2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(208, 3).Value2  = "Maule"
$ws.Cells.Item(208, 4).Value2  = 44572
$ws.Cells.Item(208, 5).Value2  = 7
$ws.Cells.Item(208, 6).Value2  = 100112043
$ws.Cells.Item(208, 7).Value2  = "Pepino ensalada"
$ws.Cells.Item(208, 8).Value2  = "Sin especificar"
$ws.Cells.Item(208, 9).Value2  = "Primera"
$ws.Cells.Item(208, 10).Value2 = 500
$ws.Cells.Item(208, 11).Value2 = 6000
$ws.Cells.Item(208, 12).Value2 = 6000
$ws.Cells.Item(208, 13).Value2 = 6000
$ws.Cells.Item(208, 14).Value2 = "`$/caja 80 unidades"
$ws.Cells.Item(208, 15).Value2 = "Región del Maule"
$ws.Cells.Item(208, 16).Value2 = 75
$ws.Cells.Item(208, 17).Value2 = 80
$ws.Cells.Item(208, 18).Value2 = "Hortaliza"
